# Fix AV reppel issue: fill in the real virement rows (row 2 had blank
# placeholders, row 3 is a brand-new beneficiary) and add a totals row (4).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$text) {
    # Assigning a plain string via .Value lets Excel auto-detect numeric-
    # looking text (e.g. long account numbers) and coerce it to a number.
    # Route it through a formula producing a text result, then collapse it
    # to a literal value via Copy/PasteSpecial (values only) so the stored
    # cell stays text (t="str"/shared string) without left-over formulas or
    # number-format styling.
    $escaped = $text.Replace("""", """""")
    $range.Formula = "=""" + $escaped + """"
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
}

# Row 2: replace blank placeholder data with the real first virement.
$ws.Range("A2").Value = "NOUBAIL MOHAMMED"
$ws.Range("B2").Value = "IR801997"
Set-TextValue $ws.Range("C2") "007400000313200019604463"
$ws.Range("D2").Value = "KHOURIBGA ZELLAKA"
$ws.Range("E2").Value = "AWB"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "004/ZZZ"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1000

# Row 3: new virement entry.
$ws.Range("A3").Value = "ZERNAKH ABDELLAH"
$ws.Range("B3").Value = "IB19558"
Set-TextValue $ws.Range("C3") "145101211406073828000084"
$ws.Range("D3").Value = "MARRAKECH BENI MELLAL"
$ws.Range("E3").Value = "BP"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "004/ZZZ"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1000

# Row 4: totals row - blank text cells with summed amounts.
$ws.Range("A4").Value = " "
$ws.Range("B4").Value = " "
$ws.Range("C4").Value = " "
$ws.Range("D4").Value = " "
$ws.Range("E4").Value = " "
$ws.Range("F4").Value = " "
$ws.Range("G4").Value = " "
$ws.Range("H4").Value = " "
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2000
